$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Database refresh -------------------------------------------------------
# Drop the oldest quarter column (old column D = "فصل دوم منتهی به 1399/06").
# Deleting the whole column shifts old E..M left into new D..L, carrying their
# values/styles/column-widths along for free.
$ws.Range("D1").EntireColumn.Delete()

# Copy column L's formatting into the (not yet populated) new column M before
# filling it in, so the brand-new cells inherit the right style/number format.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$ws.Range("M1").ColumnWidth = 30.166666666666668

# Append the newly published quarter - "فصل چهارم منتهی به 1401/12" - released
# on 1402-02-28 - as the new last column M.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-28"
$ws.Range("M11").Value = 1722165
$ws.Range("M12").Value = -2013081
$ws.Range("M13").Value = -290916
$ws.Range("M14").Value = -212458
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = -114277
$ws.Range("M17").Value = -617651
$ws.Range("M18").Value = -175917
$ws.Range("M19").Value = 2696739
$ws.Range("M20").Value = 1903171
$ws.Range("M21").Value = 147180
$ws.Range("M22").Value = 2050351
$ws.Range("M23").Value = 1713
$ws.Range("M24").Value = 2052064
$ws.Range("M25").Value = 905
$ws.Range("M26").Value = 2268000
$ws.Range("M27").Value = 905

# --- read_price algorithm change --------------------------------------------
# The "فصل چهارم منتهی به 1400/12" column (now column I after the shift above)
# was re-read under the corrected algorithm and republished later, so both its
# publish-date note and several of its figures change.
$ws.Range("I9").Value = "1402-02-28 (7)"
$ws.Range("I11").Value = 992541
$ws.Range("I13").Value = 236128
$ws.Range("I16").Value = -17178
$ws.Range("I18").Value = -122096
$ws.Range("I19").Value = 1689267
$ws.Range("I23").Value = 1117
$ws.Range("I24").Value = 1728291
